# Weekly fruit/vegetable data update: a new price-survey record for
# Alcachofa (Española variety stays, new Argentina(o) entry added) is
# inserted as row 4, pushing the existing rows 4-52 down to rows 5-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4 (shifts old rows 4..52 -> 5..53)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record
$ws.Cells.Item(4, 1).Value  = 11
$ws.Cells.Item(4, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value  = "Bíobío"
$ws.Cells.Item(4, 4).Value  = 44761
$ws.Cells.Item(4, 5).Value  = 8
$ws.Cells.Item(4, 6).Value  = 100112013
$ws.Cells.Item(4, 7).Value  = "Alcachofa"
$ws.Cells.Item(4, 8).Value  = "Argentina(o)"
$ws.Cells.Item(4, 9).Value  = "Primera"
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 13000
$ws.Cells.Item(4, 12).Value = 14000
$ws.Cells.Item(4, 13).Value = 13375
$ws.Cells.Item(4, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 268
$ws.Cells.Item(4, 17).Value = 50
$ws.Cells.Item(4, 18).Value = "Hortaliza"
